$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the "points possible" row - add Assignment 11 (column N) points
$ws.Range("N2").Value = 22

# Row 3 - Arevalo, Andres: add Assignment 7 (J), Assignment 8 (K), Assignment 9 (L)
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 11

# Row 7 - Covell, David A.: add Assignment 9 (L), Assignment 11 (N)
$ws.Range("L7").Value = 22
$ws.Range("N7").Value = 22

# Row 11 - Estrada, Andres A.: add Assignment 11 (N)
$ws.Range("N11").Value = 22

# Row 13 - Gil, Michael C.: add Assignment 2 (E), Assignment 3 (F)
$ws.Range("E13").Value = 19
$ws.Range("F13").Value = 24

# Row 15 - Gutierrez, Osvaldo: add Assignment 6 (I), Assignment 7 (J)
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 5

# Row 18 - Le, Jimmy: add Assignment 9 (L), Assignment 11 (N)
$ws.Range("L18").Value = 22
$ws.Range("N18").Value = 22

# Row 24 - Nutter, Damon A.: add Assignment 7 (J), Assignment 11 (N)
$ws.Range("J24").Value = 20
$ws.Range("N24").Value = 22

# Row 25 - Ortega Vazquez, Briana B.: add Assignment 6 (I), Assignment 11 (N)
$ws.Range("I25").Value = 20
$ws.Range("N25").Value = 22

# Row 37 - Valino, Joshua F.: add Assignment 3 (F), Assignment 4 (G)
$ws.Range("F37").Value = 24
$ws.Range("G37").Value = 13

# Row 38 - Vance, David A.: add Assignment 11 (N)
$ws.Range("N38").Value = 22

# Row 39 - Vig, Mrinal K.: add Assignment 7 (J)
$ws.Range("J39").Value = 20

# Update the active selection to reflect where the edits ended (N18)
$ws.Range("N18").Select()
